# project3_hours.xlsx - log additional tracked hours for "Duy"
# (mirrors the existing "All" row: Name | Date | Hours)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 - Duy, 3/26/2022, 1 hour
$ws.Range("A4").Value = "Duy"
$ws.Range("B4").NumberFormat = "d-mmm"
$ws.Range("B4").Value = 44646
$ws.Range("C4").Value = 1

# Row 5 - "Duy " (trailing space kept, matches source data entry), 3/30/2022, 5 hours
$ws.Range("A5").Value = "Duy "
$ws.Range("B5").NumberFormat = "d-mmm"
$ws.Range("B5").Value = 44650
$ws.Range("C5").Value = 5

# Leave the cursor where the author left it when saving
$ws.Range("C10").Select() | Out-Null
